$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.05
$ws.Range("C2").Value = 2.95
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.42
$ws.Range("F2").Value = 2.53
$ws.Range("G2").Value = 1.34
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

# Row 3
$ws.Range("C3").Value = 0.76
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.42
$ws.Range("F3").Value = 0.05
$ws.Range("G3").Value = 2.99
$ws.Range("I3").Value = 0.7

# Row 4
$ws.Range("B4").Value = 2.17
$ws.Range("C4").Value = 1.3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.42
$ws.Range("F4").Value = 2.49
$ws.Range("G4").Value = 1.32
$ws.Range("I4").Value = 0.76
